$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style cleanup around rows 137/138 -------------------------------------
# D137 currently carries the "date + underline" style (shared with nothing
# else). Move that formatting onto the new empty marker cell F138, then
# reset D137 back to the plain date format used by its neighbours (D138).

# 1) Push D137's current (date+underline) format onto F138 first, while it
#    is still intact on D137.
$ws.Range("D137").Copy()
$ws.Range("F138").PasteSpecial(-4122)   # xlPasteFormats

# 2) Reset D137 to the plain date format (copy from D138, which already
#    uses the plain date style).
$ws.Range("D138").Copy()
$ws.Range("D137").PasteSpecial(-4122)   # xlPasteFormats

# 3) F138 doesn't need a date number format (it's an empty marker cell) -
#    drop the number format back to General while keeping the underline
#    font that came across with the copy above.
$ws.Range("F138").NumberFormat = "General"

$excel.CutCopyMode = $false

# --- New "momentum" batch: rows 144-150 ------------------------------------
$newRows = @(
    @("BRFS3", 1, "momentum", 45351),
    @("CMIN3", 2, "momentum", 45351),
    @("CSMG3", 3, "momentum", 45351),
    @("ECOR3", 4, "momentum", 45351),
    @("POMO4", 5, "momentum", 45351),
    @("SBSP3", 6, "momentum", 45351),
    @("UGPA3", 7, "momentum", 45351)
)

$r = 144
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Apply the same date format used elsewhere in column D to the new rows.
$ws.Range("D2").Copy()
$ws.Range("D144:D150").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Selection / view state -------------------------------------------------
$ws.Range("D146").Select() | Out-Null
